# Add a new "time_taken" column (F) to the Peroxisomal disorders panel sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, styled like the other header cells (bold/centered/bordered).
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Per-row timestamps (metadata: how long each panel row took to process).
$timestamps = @(
    "2021-10-05 13:41:38.442006",
    "2021-10-05 13:41:38.442017",
    "2021-10-05 13:41:38.442021",
    "2021-10-05 13:41:38.442024",
    "2021-10-05 13:41:38.442027",
    "2021-10-05 13:41:38.442031",
    "2021-10-05 13:41:38.442034",
    "2021-10-05 13:41:38.442037",
    "2021-10-05 13:41:38.442040",
    "2021-10-05 13:41:38.442043",
    "2021-10-05 13:41:38.442046",
    "2021-10-05 13:41:38.442049",
    "2021-10-05 13:41:38.442052",
    "2021-10-05 13:41:38.442055",
    "2021-10-05 13:41:38.442058",
    "2021-10-05 13:41:38.442061",
    "2021-10-05 13:41:38.442064",
    "2021-10-05 13:41:38.442068",
    "2021-10-05 13:41:38.442071",
    "2021-10-05 13:41:38.442074",
    "2021-10-05 13:41:38.442077",
    "2021-10-05 13:41:38.442080",
    "2021-10-05 13:41:38.442083",
    "2021-10-05 13:41:38.442086",
    "2021-10-05 13:41:38.442089",
    "2021-10-05 13:41:38.442092",
    "2021-10-05 13:41:38.442095",
    "2021-10-05 13:41:38.442098",
    "2021-10-05 13:41:38.442102",
    "2021-10-05 13:41:38.442105",
    "2021-10-05 13:41:38.442108",
    "2021-10-05 13:41:38.442111",
    "2021-10-05 13:41:38.442114",
    "2021-10-05 13:41:38.442118",
    "2021-10-05 13:41:38.442121",
    "2021-10-05 13:41:38.442124",
    "2021-10-05 13:41:38.442127",
    "2021-10-05 13:41:38.442130"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
